$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.818.91'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '2.238.08'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.91'
$ws.Range("E5").Value = '  +2.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.56'
$ws.Range("E6").Value = '  +8.65%  '
$ws.Range("E7").Value = '  -2.94%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.56'
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0929'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.18'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("E13").Value = '  -2.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.881'
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").Value = '2.574.77'
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").Value = '2.252.41'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '42.761.29'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000108'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.64'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("B22").Value = 'PancakeSwap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.22'
$ws.Range("E22").Value = '  +11.81%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.36'
$ws.Range("E23").Value = '  -2.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '231.71'
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.15'
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.03'
$ws.Range("E26").Value = '  +5.96%  '
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.95'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.27'
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.30'
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.24'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.51'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.16'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0902'
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.61'
$ws.Range("E35").Value = '  +18.86%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.58'
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.128'
$ws.Range("E37").Value = '  -2.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0373'
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.65'
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.18'
$ws.Range("E42").Value = '  -3.63%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.45'
$ws.Range("E43").Value = '  -5.34%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.234'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.33'
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.58'
$ws.Range("E47").Value = '  -8.17%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.30'
$ws.Range("E48").Value = '  +3.08%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.54'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.650'
$ws.Range("E50").Value = '  +8.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '101.14'
$ws.Range("E51").Value = '  +1.14%  '
